# The workbook is already open; grab the workbook and its (only) active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" -> "Pricings"
$ws.Name = "Pricings"

# Move the sheet's selection from E4 to D22
$ws.Range("D22").Select()
